# "M11 Froze Encoder 1234" - refresh the per-epoch accuracy numbers pulled
# from the latest training run of the (re-started / "frozen encoder") model,
# and re-stamp the stale Python object repr left in column A for the tail
# rows that came from the notebook's display-object placeholder.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New accuracy readings (column B) for the affected epochs (row = sheet row).
$newAccuracy = @{
    3   = 0.828125
    4   = 0.796875
    5   = 0.8125
    6   = 0.8125
    7   = 0.71875
    9   = 0.71875
    11  = 0.65625
    12  = 0.640625
    14  = 0.6875
    15  = 0.671875
    16  = 0.671875
    18  = 0.59375
    19  = 0.59375
    20  = 0.59375
    21  = 0.578125
    22  = 0.5
    23  = 0.5
    24  = 0.5
    25  = 0.5
    26  = 0.5
    27  = 0.5
    28  = 0.5
    29  = 0.5
    30  = 0.5
    31  = 0.5
    32  = 0.5
    33  = 0.5
    34  = 0.5
    35  = 0.5
    36  = 0.5
    37  = 0.5
    38  = 0.5
    39  = 0.5
    40  = 0.5
    41  = 0.5
    42  = 0.5
    43  = 0.5
    44  = 0.5
    45  = 0.5
    46  = 0.5
    47  = 0.5
    48  = 0.5
    49  = 0.5
    50  = 0.5
    51  = 0.5
    52  = 0.5
    53  = 0.5
    54  = 0.5
    55  = 0.5
    56  = 0.5
    57  = 0.5
    58  = 0.5
    59  = 0.5
    60  = 0.5
    61  = 0.5
    62  = 0.5
    63  = 0.5
    64  = 0.5
    65  = 0.5
    66  = 0.5
    67  = 0.5
    68  = 0.5
    69  = 0.5
    70  = 0.5
    71  = 0.5
    72  = 0.5
    73  = 0.5
    74  = 0.5
    75  = 0.5
    76  = 0.5
    77  = 0.5
    78  = 0.5
    79  = 0.5
    80  = 0.5
    81  = 0.5
    82  = 0.5
    83  = 0.5
    84  = 0.5
    85  = 0.5
    86  = 0.5
    87  = 0.5
    88  = 0.5
    89  = 0.5
    90  = 0.5
    91  = 0.5
    92  = 0.5
    93  = 0.5
    94  = 0.515625
    95  = 0.515625
    96  = 0.515625
    97  = 0.515625
    98  = 0.515625
    99  = 0.515625
    100 = 0.515625
    101 = 0.515625
    102 = 0.515625
    103 = 0.5625
    104 = 0.578125
    105 = 0.578125
    106 = 0.65625
    107 = 0.625
    108 = 0.546875
    109 = 0.515625
    110 = 0.625
    112 = 0.53125
    113 = 0.65625
    114 = 0.484375
    116 = 0.6764705882352942
}

foreach ($row in $newAccuracy.Keys) {
    $ws.Cells.Item($row, 2).Value = $newAccuracy[$row]
}

# The tail rows' column A held the repr() of a live Python object
# (`<__main__.DisplayOutputs object at 0x...>`) captured when the notebook
# cell was last executed; re-running it changed the object's memory address.
$oldRepr = "<__main__.DisplayOutputs object at 0x7f04f014d820>"
$newRepr = "<__main__.DisplayOutputs object at 0x7f865009d460>"
for ($row = 102; $row -le 116; $row++) {
    $cell = $ws.Cells.Item($row, 1)
    if ($cell.Value2 -eq $oldRepr) {
        $cell.Value = $newRepr
    }
}

# Select the whole sheet (as left behind by the notebook export step) instead
# of just the populated A2:B116 block.
$ws.Cells.Select()
